$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting (bold, centered, bordered) used by the other header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in team record values (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 30).Value = 53
    $ws.Cells.Item($r, 31).Value = 60
    $ws.Cells.Item($r, 32).Value = 0
}
